$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # The "Price" column holds values that look numeric (e.g. "0.9995",
    # "1.000", "29.186.14") but must stay plain text, exactly as in the
    # source data. Force text format while assigning, then restore the
    # default "Normal" style so no stray style index is left behind.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "29.186.14"
$ws.Range("E2").Value = "  +0.31%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.843.06"
$ws.Range("E3").Value = "  +0.62%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "0.9995"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5 - BNB
$ws.Range("E5").Value = "  -0.06%  "

# Row 6 - XRP
Set-TextValue "D6" "0.6709"
$ws.Range("E6").Value = "  -1.95%  "

# Row 7 - USDC
Set-TextValue "D7" "1.000"
$ws.Range("E7").Value = "  -0.01%  "

# Row 8 - Dogecoin
Set-TextValue "D8" "0.07419"
$ws.Range("E8").Value = "  -0.77%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.2953"
$ws.Range("E9").Value = "  -2.11%  "

# Row 10 - Solana
Set-TextValue "D10" "22.83"
$ws.Range("E10").Value = "  -1.25%  "

# Row 11 - TRON
Set-TextValue "D11" "0.07715"
$ws.Range("E11").Value = "  +0.75%  "

# Row 12 - WrappedEther
Set-TextValue "D12" "1.835.08"
$ws.Range("E12").Value = "  -0.06%  "

# Row 13 - Polkadot
Set-TextValue "D13" "5.003"
$ws.Range("E13").Value = "  -1.25%  "

# Row 14 - Polygon
Set-TextValue "D14" "0.6763"
$ws.Range("E14").Value = "  -0.98%  "

# Row 15 - Litecoin
Set-TextValue "D15" "86.03"
$ws.Range("E15").Value = "  -1.87%  "

# Row 16 - Uniswap
Set-TextValue "D16" "6.131"
$ws.Range("E16").Value = "  -0.60%  "

# Row 17 - now ShibaInu (was WrappedBTC)
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D17" "0.000008303"
$ws.Range("E17").Value = "  +1.48%  "

# Row 18 - now WrappedBTC (was ShibaInu)
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D18" "29.113.45"
$ws.Range("E18").Value = "  +0.01%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "228.49"
$ws.Range("E19").Value = "  +0.65%  "

# Row 20 - Avalanche
Set-TextValue "D20" "12.53"
$ws.Range("E20").Value = "  -0.11%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.10%  "

# Row 22 - Chainlink
Set-TextValue "D22" "7.192"
$ws.Range("E22").Value = "  -3.18%  "

# Row 23 - BinanceUSD
Set-TextValue "D23" "1.001"
$ws.Range("E23").Value = "  +0.00%  "

# Row 24 - Monero
Set-TextValue "D24" "160.50"
$ws.Range("E24").Value = "  +0.18%  "

# Row 25 - Cosmos
Set-TextValue "D25" "8.696"
$ws.Range("E25").Value = "  -0.87%  "

# Row 26 - Stellar
Set-TextValue "D26" "0.1404"
$ws.Range("E26").Value = "  -3.91%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "18.00"
$ws.Range("E27").Value = "  -0.69%  "

# Row 28 - PancakeSwap
Set-TextValue "D28" "1.508"
$ws.Range("E28").Value = "  -0.26%  "

# Row 29 - Filecoin
Set-TextValue "D29" "4.169"
$ws.Range("E29").Value = "  -2.47%  "

# Row 30 - InternetComputer(DFINITY)
Set-TextValue "D30" "4.070"
$ws.Range("E30").Value = "  -1.88%  "

# Row 31 - Toncoin
Set-TextValue "D31" "1.189"
$ws.Range("E31").Value = "  -0.62%  "

# Row 32 - Hedera
Set-TextValue "D32" "0.05321"
$ws.Range("E32").Value = "  +2.95%  "

# Row 33 - now ImmutableX (was LidoDAOToken)
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D33" "0.7602"
$ws.Range("E33").Value = "  -1.15%  "

# Row 34 - now LidoDAOToken (was ImmutableX)
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D34" "1.880"
$ws.Range("E34").Value = "  +1.79%  "

# Row 35 - ARBITRUM
Set-TextValue "D35" "1.137"
$ws.Range("E35").Value = "  +0.43%  "

# Row 36 - HuobiToken
Set-TextValue "D36" "2.686"
$ws.Range("E36").Value = "  +0.47%  "

# Row 37 - Maker
Set-TextValue "D37" "1.332.00"
$ws.Range("E37").Value = "  +1.65%  "

# Row 38 - VeChain
Set-TextValue "D38" "0.01799"
$ws.Range("E38").Value = "  -2.15%  "

# Row 39 - MXToken
Set-TextValue "D39" "2.736"
$ws.Range("E39").Value = "  +0.39%  "

# Row 40 - TrustWalletToken
Set-TextValue "D40" "0.9249"
$ws.Range("E40").Value = "  -0.96%  "

# Row 41 - FraxShare
Set-TextValue "D41" "5.968"
$ws.Range("E41").Value = "  +3.03%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.19%  "

# Row 43 - Quant
Set-TextValue "D43" "103.42"
$ws.Range("E43").Value = "  -1.02%  "

# Row 44 - XinFinNetwork
Set-TextValue "D44" "0.07990"
$ws.Range("E44").Value = "  +12.51%  "

# Row 45 - BabyDogeCoin
Set-TextValue "D45" "0.00000000123"
$ws.Range("E45").Value = "  -0.29%  "

# Row 46 - RocketPoolETH
Set-TextValue "D46" "1.968.95"
$ws.Range("E46").Value = "  -0.61%  "

# Row 47 - Mantle
Set-TextValue "D47" "0.5163"
$ws.Range("E47").Value = "  -0.53%  "

# Row 48 - RenderToken
$ws.Range("E48").Value = "  -0.36%  "

# Row 49 - Aave
Set-TextValue "D49" "63.95"
$ws.Range("E49").Value = "  -2.18%  "

# Row 50 - EnergySwap
Set-TextValue "D50" "9.187"
$ws.Range("E50").Value = "  -4.44%  "

# Row 51 - Cronos
Set-TextValue "D51" "0.05943"
$ws.Range("E51").Value = "  +0.22%  "
